$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D whose text values look numeric need an explicit
# Text format while being written so Excel does not coerce them to
# numbers (stripping trailing zeros / using scientific notation).
# Formats are cleared again immediately after so no extra style is
# left behind on the cell.

$ws.Range("D2").Value = "62.952.59"
$ws.Range("E2").Value = "  +2.73%  "

$ws.Range("D3").Value = "3.025.37"
$ws.Range("E3").Value = "  +1.56%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  -0.07%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "595.30"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.29%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "152.92"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +6.77%  "

$ws.Range("E7").Value = "  +0.05%  "

$ws.Range("D8").Value = "3.022.72"
$ws.Range("E8").Value = "  +1.61%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.513"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -0.01%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.98"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +15.81%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.150"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +1.88%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.463"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +2.39%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000233"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +2.96%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "35.64"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +4.54%  "

$ws.Range("E15").Value = "  -1.33%  "

$ws.Range("D16").Value = "3.526.78"
$ws.Range("E16").Value = "  +1.51%  "

$ws.Range("B17").Value = "Polkadot"
$ws.Range("C17").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "7.07"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +2.51%  "

$ws.Range("B18").Value = "WrappedBTC"
$ws.Range("C18").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D18").Value = "62.898.15"
$ws.Range("E18").Value = "  +2.72%  "

$ws.Range("D19").Value = "3.023.78"
$ws.Range("E19").Value = "  +1.59%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "449.40"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +0.48%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.22"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +1.62%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.696"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +2.10%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.53"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +2.65%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "82.97"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +1.63%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "11.42"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +8.44%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.31"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +6.31%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.35"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +3.00%  "

$ws.Range("E28").Value = "  +0.00%  "

$ws.Range("E29").Value = "  +4.78%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.28"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +11.24%  "

$ws.Range("E31").Value = "  +0.56%  "

$ws.Range("E32").Value = "  -0.12%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "27.66"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +2.08%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.111"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +1.09%  "

$ws.Range("D35").Value = "0.0₃0873"
$ws.Range("E35").Value = "  +7.32%  "

$ws.Range("E36").Value = "  +2.89%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.88"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +2.06%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.15"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +10.74%  "

$ws.Range("B39").Value = "Kaspa"
$ws.Range("C39").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.130"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +8.38%  "

$ws.Range("B40").Value = "Stacks"
$ws.Range("C40").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.10"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +3.53%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "50.52"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +0.84%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "9.03"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +0.38%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "44.48"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +16.68%  "

$ws.Range("E44").Value = "  +13.95%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "391.45"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +1.21%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0359"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +2.43%  "

$ws.Range("D47").Value = "2.711.40"
$ws.Range("E47").Value = "  +1.01%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "133.76"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +2.84%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "26.73"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +14.52%  "

$ws.Range("E50").Value = "  -0.03%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.28"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +6.21%  "
